# Daily attendance processing - 2026-01-18 15:58:31
# Reorder the "Recorded By" (column G) entries so that "System" is listed
# first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# This only applies to entries that are exactly "<single address>, System"
# (entries that already start with "System" or that contain the
# "backup@backdoor.com" account stay untouched, matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = 7
    $value = $cell.Text

    $parts = $value -split ', '
    if ($parts.Count -eq 2 -and $parts[1] -eq 'System' -and $parts[0] -ne 'System' -and $parts[0] -ne 'backup@backdoor.com') {
        $cell.Value = "System, " + $parts[0]
    }
}
